$wb = $excel.ActiveWorkbook

# ALC!row5
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 333.5625
$ws.Range("I5").Value = 188.14285
$ws.Range("J5").Value = 446.66666
$ws.Range("K5").Value = 188.14285
$ws.Range("L5").Value = 446.66666
$ws.Range("M5").Value = -73.14285000000001
$ws.Range("N5").Value = -676.66666

# ALC!row15
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 492.76135
$ws.Range("I15").Value = 492.76135
$ws.Range("K15").Value = 1478.28405
$ws.Range("M15").Value = -1309.28405

# ALC!row19
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1097431.5
$ws.Range("I19").Value = 1880321.1
$ws.Range("J19").Value = 1386.2
$ws.Range("K19").Value = 1880321.1
$ws.Range("L19").Value = 1386.2
$ws.Range("M19").Value = -1880146.1
$ws.Range("N19").Value = -1736.2

# ALC!row29
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 875
$ws.Range("I29").Value = 875
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 2625
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -2344
$ws.Range("N29").ClearContents()

# ALC!row32
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1209
$ws.Range("J32").Value = 1342.5714
$ws.Range("L32").Value = 1342.5714
$ws.Range("N32").Value = -1994.5714

# ALC!row33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 176.20833
$ws.Range("I33").Value = 161.19048
$ws.Range("K33").Value = 161.19048
$ws.Range("M33").Value = 67.80951999999999

# ALC!row55
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 681.875
$ws.Range("I55").Value = 565.5
$ws.Range("J55").Value = 720.6667
$ws.Range("K55").Value = 565.5
$ws.Range("L55").Value = 720.6667
$ws.Range("M55").Value = -351.5
$ws.Range("N55").Value = -1148.6667

# ALC!row62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 1359.8889
$ws.Range("I62").Value = 1359.8889
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 1359.8889
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -735.8888999999999
$ws.Range("N62").ClearContents()

# ALC!row65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 1359.8889
$ws.Range("I65").Value = 1359.8889
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 6799.4445
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -3679.4445
$ws.Range("N65").ClearContents()

# ALC!row93
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H93").Value = 33407.145
$ws.Range("J93").Value = 33407.145
$ws.Range("L93").Value = 33407.145
$ws.Range("N93").Value = -38399.145

# ALC!row112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1269.0625
$ws.Range("J112").Value = 1269.0625
$ws.Range("L112").Value = 3807.1875
$ws.Range("N112").Value = -6023.1875

# ALC!row113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 8511
$ws.Range("I113").Value = 2958.75
$ws.Range("J113").Value = 9991.6
$ws.Range("K113").Value = 2958.75
$ws.Range("L113").Value = 9991.6
$ws.Range("M113").Value = 295.25
$ws.Range("N113").Value = -16499.6

# ALC!row123
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H123").Value = 41920
$ws.Range("J123").Value = 41920
$ws.Range("L123").Value = 41920
$ws.Range("N123").Value = -51720

# ALC!row132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 327098.62
$ws.Range("I132").Value = 103795.766
$ws.Range("J132").Value = 2504301.5
$ws.Range("K132").Value = 311387.298
$ws.Range("L132").Value = 7512904.5
$ws.Range("M132").Value = -308857.298
$ws.Range("N132").Value = -7517964.5

# ALC!row137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1906848.4
$ws.Range("I137").Value = 4329977
$ws.Range("J137").Value = 2961.5
$ws.Range("K137").Value = 12989931
$ws.Range("L137").Value = 8884.5
$ws.Range("M137").Value = -12987381
$ws.Range("N137").Value = -13984.5

# ALC!row138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2667.2368
$ws.Range("I138").Value = 1830.7693
$ws.Range("J138").Value = 3102.2
$ws.Range("K138").Value = 5492.3079
$ws.Range("L138").Value = 9306.599999999999
$ws.Range("M138").Value = -352.3078999999998
$ws.Range("N138").Value = -19586.6

# ARM!row32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5615.5283
$ws.Range("I32").Value = 5929.091
$ws.Range("J32").Value = 5098.15
$ws.Range("K32").Value = 5929.091
$ws.Range("L32").Value = 5098.15
$ws.Range("M32").Value = -5642.091
$ws.Range("N32").Value = -5672.15

# ARM!row64
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

# ARM!row67
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

# ARM!row101
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H101").Value = 37120.4
$ws.Range("J101").Value = 37120.4
$ws.Range("L101").Value = 37120.4
$ws.Range("N101").Value = -43610.4

# ARM!row132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2150.257
$ws.Range("I132").Value = 1380.7916
$ws.Range("J132").Value = 3829.0908
$ws.Range("K132").Value = 4142.3748
$ws.Range("L132").Value = 11487.2724
$ws.Range("M132").Value = -1612.3748
$ws.Range("N132").Value = -16547.2724

# BSM!row62
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 37999.5
$ws.Range("J62").Value = 37999.5
$ws.Range("L62").Value = 37999.5
$ws.Range("N62").Value = -39371.5

# BSM!row65
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H65").Value = 37999.5
$ws.Range("J65").Value = 37999.5
$ws.Range("L65").Value = 113998.5
$ws.Range("N65").Value = -120862.5

# BSM!row82
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 23700.75
$ws.Range("I82").Value = 3720
$ws.Range("J82").Value = 33691.125
$ws.Range("K82").Value = 3720
$ws.Range("L82").Value = 33691.125
$ws.Range("M82").Value = -3337
$ws.Range("N82").Value = -34457.125

# BSM!row85
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H85").Value = 23700.75
$ws.Range("I85").Value = 3720
$ws.Range("J85").Value = 33691.125
$ws.Range("K85").Value = 3720
$ws.Range("L85").Value = 33691.125
$ws.Range("M85").Value = -2394
$ws.Range("N85").Value = -36343.125

# BSM!row94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 825
$ws.Range("I94").Value = 766.6667
$ws.Range("J94").Value = 1000
$ws.Range("K94").Value = 766.6667
$ws.Range("L94").Value = 1000
$ws.Range("M94").Value = -315.6667
$ws.Range("N94").Value = -1902

# BSM!row134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2193.7874
$ws.Range("I134").Value = 1335.579
$ws.Range("J134").Value = 5817.3335
$ws.Range("K134").Value = 4006.737
$ws.Range("L134").Value = 17452.0005
$ws.Range("M134").Value = -1471.737
$ws.Range("N134").Value = -22522.0005

# CRP!row31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7492.244
$ws.Range("I31").Value = 3297.2144
$ws.Range("J31").Value = 9667.444
$ws.Range("K31").Value = 3297.2144
$ws.Range("L31").Value = 9667.444
$ws.Range("M31").Value = -3002.2144
$ws.Range("N31").Value = -10257.444

# CRP!row34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 7492.244
$ws.Range("I34").Value = 3297.2144
$ws.Range("J34").Value = 9667.444
$ws.Range("K34").Value = 3297.2144
$ws.Range("L34").Value = 9667.444
$ws.Range("M34").Value = -3095.2144
$ws.Range("N34").Value = -10071.444

# CRP!row54
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H54").Value = 11899
$ws.Range("J54").Value = 11899
$ws.Range("L54").Value = 11899
$ws.Range("N54").Value = -13215

# CRP!row58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2174.42
$ws.Range("I58").Value = 1372.8611
$ws.Range("K58").Value = 1372.8611
$ws.Range("M58").Value = -1169.8611

# CRP!row98
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H98").Value = 39141
$ws.Range("J98").Value = 39141
$ws.Range("L98").Value = 39141
$ws.Range("N98").Value = -43633

# CRP!row132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2056.1943
$ws.Range("I132").Value = 1633.6666
$ws.Range("J132").Value = 4168.8335
$ws.Range("K132").Value = 4900.9998
$ws.Range("L132").Value = 12506.5005
$ws.Range("M132").Value = -2370.9998
$ws.Range("N132").Value = -17566.5005

# CRP!row134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 4526.2856
$ws.Range("I134").Value = 4830.222
$ws.Range("K134").Value = 14490.666
$ws.Range("M134").Value = -11955.666

# CRP!row136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2174.42
$ws.Range("I136").Value = 1372.8611
$ws.Range("K136").Value = 4118.5833
$ws.Range("M136").Value = -1568.5833

# CUL!row23
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 139.76923
$ws.Range("I23").Value = 53.5
$ws.Range("J23").Value = 178.11111
$ws.Range("K23").Value = 160.5
$ws.Range("L23").Value = 534.3333299999999
$ws.Range("M23").Value = 74.5
$ws.Range("N23").Value = -1004.33333

# CUL!row104
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 9000
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 9000
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 27000
$ws.Range("M104").ClearContents()
$ws.Range("N104").Value = -32242

# CUL!row131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 807.28
$ws.Range("J131").Value = 832.92633
$ws.Range("L131").Value = 2498.77899
$ws.Range("N131").Value = -12578.77899

# CUL!row132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2250.4814
$ws.Range("I132").Value = 998.63635
$ws.Range("J132").Value = 3111.125
$ws.Range("K132").Value = 8987.727150000001
$ws.Range("L132").Value = 28000.125
$ws.Range("M132").Value = -6457.727150000001
$ws.Range("N132").Value = -33060.125

# GSM!row113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1671.6428
$ws.Range("I113").Value = 1553
$ws.Range("J113").Value = 1885.2
$ws.Range("K113").Value = 1553
$ws.Range("L113").Value = 1885.2
$ws.Range("M113").Value = 617
$ws.Range("N113").Value = -6225.2

# LTW!row110
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H110").Value = 27883.857
$ws.Range("J110").Value = 27883.857
$ws.Range("L110").Value = 27883.857
$ws.Range("N110").Value = -36063.857

# LTW!row132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5042.0415
$ws.Range("I132").Value = 4137.625
$ws.Range("J132").Value = 6850.875
$ws.Range("K132").Value = 12412.875
$ws.Range("L132").Value = 20552.625
$ws.Range("M132").Value = -9882.875
$ws.Range("N132").Value = -25612.625

# WVR!row40
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 13654.454
$ws.Range("I40").Value = 4000
$ws.Range("J40").Value = 14619.9
$ws.Range("K40").Value = 4000
$ws.Range("L40").Value = 14619.9
$ws.Range("M40").Value = -3851
$ws.Range("N40").Value = -14917.9

# WVR!row107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 526.9091
$ws.Range("I107").Value = 283.2857
$ws.Range("J107").Value = 953.25
$ws.Range("K107").Value = 849.8571000000001
$ws.Range("L107").Value = 2859.75
$ws.Range("M107").Value = 1070.1429
$ws.Range("N107").Value = -6699.75

# WVR!row116
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H116").Value = 34000
$ws.Range("J116").Value = 34000
$ws.Range("L116").Value = 34000
$ws.Range("N116").Value = -43178

# WVR!row126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 508922.66
$ws.Range("I126").Value = 1073.375
$ws.Range("J126").Value = 2134040.5
$ws.Range("K126").Value = 3220.125
$ws.Range("L126").Value = 6402121.5
$ws.Range("M126").Value = -750.125
$ws.Range("N126").Value = -6407061.5

# WVR!row136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3318.8064
$ws.Range("I136").Value = 994.4545000000001
$ws.Range("J136").Value = 9000.556
$ws.Range("K136").Value = 2983.3635
$ws.Range("L136").Value = 27001.668
$ws.Range("M136").Value = -433.3635000000004
$ws.Range("N136").Value = -32101.668
